# "Added Helper Package and Classes"
#
# The class-diagram style table on Tabelle1 lists, for every class block,
# the class name + its methods, followed by a bold "author" row.
# The "Initializer / getNameFromUser()" class block (row 40) is removed,
# and the now-orphaned author cell that used to sit below it (row 41,
# which held "LexuTros") is cleared too, while keeping its bold styling.
# The following class blocks (DiceRoller/rollDice, Mondra) keep their
# original row numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A40").Value = ""
$ws.Range("B40").Value = ""
$ws.Range("A41").Value = ""

# Leave the selection where the user ended their edit.
$ws.Range("A40").Select()
